$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# A leading apostrophe forces Excel to treat the value as literal text,
# preserving formats such as "322.30", "1.00", multi-dot prices, and
# the padded "  +1.23%  " volume strings exactly as authored.

$ws.Range("D2").Value = "'42.460.89"
$ws.Range("E2").Value = "'  +0.87%  "
$ws.Range("D3").Value = "'2.286.52"
$ws.Range("E3").Value = "'  -0.86%  "
$ws.Range("E4").Value = "'  +0.27%  "
$ws.Range("D5").Value = "'322.30"
$ws.Range("E5").Value = "'  +1.70%  "
$ws.Range("D6").Value = "'102.74"
$ws.Range("E6").Value = "'  -2.30%  "
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E8").Value = "'  +0.25%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "'  -0.55%  "
$ws.Range("D10").Value = "'39.93"
$ws.Range("E10").Value = "'  +0.20%  "
$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "'  -0.15%  "
$ws.Range("D12").Value = "'8.33"
$ws.Range("E12").Value = "'  -1.79%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "'  -0.40%  "
$ws.Range("D14").Value = "'0.967"
$ws.Range("E14").Value = "'  -1.42%  "
$ws.Range("D15").Value = "'15.16"
$ws.Range("E15").Value = "'  -1.98%  "
$ws.Range("D16").Value = "'2.632.63"
$ws.Range("D17").Value = "'2.284.13"
$ws.Range("E17").Value = "'  -1.23%  "
$ws.Range("D18").Value = "'42.398.45"
$ws.Range("E18").Value = "'  +0.71%  "
$ws.Range("D19").Value = "'7.39"
$ws.Range("E19").Value = "'  -4.60%  "
$ws.Range("E20").Value = "'  -0.58%  "
$ws.Range("D21").Value = "'12.86"
$ws.Range("E21").Value = "'  +28.44%  "
$ws.Range("E22").Value = "'  +1.71%  "
$ws.Range("D23").Value = "'73.06"
$ws.Range("D24").Value = "'268.35"
$ws.Range("E24").Value = "'  -5.31%  "
$ws.Range("E25").Value = "'  -3.22%  "
$ws.Range("E26").Value = "'  -0.16%  "
$ws.Range("D27").Value = "'10.87"
$ws.Range("E27").Value = "'  -1.02%  "
$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "'  +4.15%  "
$ws.Range("D29").Value = "'22.44"
$ws.Range("E29").Value = "'  -3.79%  "
$ws.Range("D30").Value = "'37.98"
$ws.Range("E30").Value = "'  +6.68%  "
$ws.Range("D31").Value = "'164.37"
$ws.Range("E31").Value = "'  -0.80%  "
$ws.Range("D32").Value = "'6.12"
$ws.Range("E32").Value = "'  +3.27%  "
$ws.Range("D33").Value = "'0.0879"
$ws.Range("E33").Value = "'  -0.62%  "
$ws.Range("E34").Value = "'  +0.54%  "
$ws.Range("E35").Value = "'  -13.03%  "
$ws.Range("E36").Value = "'  -3.64%  "
$ws.Range("E37").Value = "'  -1.76%  "
$ws.Range("D38").Value = "'0.0354"
$ws.Range("E38").Value = "'  +0.17%  "
$ws.Range("E39").Value = "'  +1.59%  "
$ws.Range("E40").Value = "'  -5.95%  "
$ws.Range("D41").Value = "'1.52"
$ws.Range("E41").Value = "'  +1.59%  "
$ws.Range("D42").Value = "'68.96"
$ws.Range("E42").Value = "'  -2.68%  "
$ws.Range("B43").Value = "'FirstDigitalUSD"
$ws.Range("C43").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "'  +0.22%  "
$ws.Range("B44").Value = "'Algorand"
$ws.Range("C44").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.225"
$ws.Range("E44").Value = "'  -1.18%  "
$ws.Range("B45").Value = "'Celestia"
$ws.Range("C45").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'12.38"
$ws.Range("E45").Value = "'  +1.74%  "
$ws.Range("B46").Value = "'BitcoinSV"
$ws.Range("C46").Value = "'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'90.92"
$ws.Range("E46").Value = "'  -11.30%  "
$ws.Range("D47").Value = "'113.40"
$ws.Range("E47").Value = "'  -3.39%  "
$ws.Range("D48").Value = "'80.09"
$ws.Range("E48").Value = "'  +1.47%  "
$ws.Range("D49").Value = "'8.92"
$ws.Range("E49").Value = "'  -2.75%  "
$ws.Range("D50").Value = "'5.23"
$ws.Range("E50").Value = "'  -2.12%  "
$ws.Range("D51").Value = "'1.593.49"
$ws.Range("E51").Value = "'  +1.96%  "
